$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.08873594589893813
$ws.Range("C2").Value = 0.5038494199792003
$ws.Range("D2").Value = -0.0804002763435806
$ws.Range("E2").Value = 0.7915875151041314
$ws.Range("F2").Value = 0.7713506143760337
$ws.Range("G2").Value = 0.3749201162859844
$ws.Range("H2").Value = 0.5168869453501342
$ws.Range("I2").Value = 0.6892623450793038
$ws.Range("J2").Value = 0.1844339823288103
$ws.Range("K2").Value = 0.4621062954844631

# Row 3
$ws.Range("B3").Value = 0.526277399612209
$ws.Range("C3").Value = -0.07779673678113191
$ws.Range("D3").Value = 0.7695198294501161
$ws.Range("E3").Value = 0.7653686340716788
$ws.Range("F3").Value = 0.3695768661035733
$ws.Range("G3").Value = 0.506518411979768
$ws.Range("H3").Value = 0.68078502295568
$ws.Range("I3").Value = 0.1766616740219995
$ws.Range("J3").Value = 0.4534823219514945
$ws.Range("K3").Value = 0.2839431369332225

# Row 4
$ws.Range("B4").Value = -0.05637216532391182
$ws.Range("C4").Value = 0.8896976782493284
$ws.Range("D4").Value = 0.6761742896578956
$ws.Range("E4").Value = 0.3398498622549955
$ws.Range("F4").Value = 0.5151445320096781
$ws.Range("G4").Value = 0.6548448495302448
$ws.Range("H4").Value = 0.153883110993772
$ws.Range("I4").Value = 0.4398642868028766
$ws.Range("J4").Value = 0.2654223397480467
$ws.Range("K4").Value = 0.570669944985061

# Row 5
$ws.Range("B5").Value = 0.847377045928939
$ws.Range("C5").Value = 0.6469698158021624
$ws.Range("D5").Value = 0.3368292624500743
$ws.Range("E5").Value = 0.4962832483981977
$ws.Range("F5").Value = 0.6345137184650405
$ws.Range("G5").Value = 0.1393529555595242
$ws.Range("H5").Value = 0.4232232413106087
$ws.Range("I5").Value = 0.2478384943192965
$ws.Range("J5").Value = 0.5541960614550182
$ws.Range("K5").Value = -0.05728328644410208

# Row 6
$ws.Range("B6").Value = 0.9871675564200725
$ws.Range("C6").Value = 0.4122003242340114
$ws.Range("D6").Value = 0.3051899620851986
$ws.Range("E6").Value = 0.6592076310517737
$ws.Range("F6").Value = 0.1487323591158202
$ws.Range("G6").Value = 0.3684124426992176
$ws.Range("H6").Value = 0.2289444034306267
$ws.Range("I6").Value = 0.5393320606399725
$ws.Range("J6").Value = -0.0859949970734728
$ws.Range("K6").Value = 0.6071339948549791

# Row 7
$ws.Range("B7").Value = 0.862895196224262
$ws.Range("C7").Value = 0.352300664297557
$ws.Range("D7").Value = 0.4189244002609654
$ws.Range("E7").Value = 0.184677440181683
$ws.Range("F7").Value = 0.3663616852596248
$ws.Range("G7").Value = 0.1523761639945965
$ws.Range("H7").Value = 0.511932666264689
$ws.Range("I7").Value = -0.1122840472711982
$ws.Range("J7").Value = 0.5637367041416466
$ws.Range("K7").Value = 0.2970525035592049

# Row 8
$ws.Range("B8").Value = 0.6646262512210954
$ws.Range("C8").Value = 0.5522131399964898
$ws.Range("D8").Value = 0.005599018365491398
$ws.Range("E8").Value = 0.394555924030192
$ws.Range("F8").Value = 0.18840147518699
$ws.Range("G8").Value = 0.4743910960604755
$ws.Range("H8").Value = -0.1178263863585594
$ws.Range("I8").Value = 0.5679120330803951
$ws.Range("J8").Value = 0.2858677898194339

# Row 9
$ws.Range("B9").Value = 0.787803631104331
$ws.Range("C9").Value = 0.09027775923980097
$ws.Range("D9").Value = 0.2489555573964748
$ws.Range("E9").Value = 0.1983552180462326
$ws.Range("F9").Value = 0.4897562657600204
$ws.Range("G9").Value = -0.1566747213159825
$ws.Range("H9").Value = 0.5520662240532093
$ws.Range("I9").Value = 0.2775335613519331

# Row 10
$ws.Range("B10").Value = 0.4013019457211838
$ws.Range("C10").Value = 0.36604433180767
$ws.Range("D10").Value = 0.03589107659666579
$ws.Range("E10").Value = 0.518715216225222
$ws.Range("F10").Value = -0.1209318488610789
$ws.Range("G10").Value = 0.5207385776695821
$ws.Range("H10").Value = 0.2743085116504074

# Row 11
$ws.Range("B11").Value = 0.6128695092117844
$ws.Range("C11").Value = 0.0535469441345553
$ws.Range("D11").Value = 0.4240929771142275
$ws.Range("E11").Value = -0.0887144606125988
$ws.Range("F11").Value = 0.5331267034972994
$ws.Range("G11").Value = 0.2534447081011285

# Row 12
$ws.Range("B12").Value = 0.2932233035507672
$ws.Range("C12").Value = 0.509117220583441
$ws.Range("D12").Value = -0.2051990389706129
$ws.Range("E12").Value = 0.5482319751491519
$ws.Range("F12").Value = 0.2766837437271186

# Row 13
$ws.Range("B13").Value = 0.6739772976175282
$ws.Range("C13").Value = -0.1916617667226967
$ws.Range("D13").Value = 0.4852590561591889
$ws.Range("E13").Value = 0.2867219094086165

# Row 14
$ws.Range("B14").Value = 0.06218727514271133
$ws.Range("C14").Value = 0.5845771063412253
$ws.Range("D14").Value = 0.1751453671933744

# Row 15
$ws.Range("B15").Value = 0.6286367975806744
$ws.Range("C15").Value = 0.1965658720679752

# Row 16
$ws.Range("B16").Value = 0.4328090033804217
